# Fruta / hortaliza, semanal
# The weekly refresh re-sorted the "Frambuesa" (Vega Central Mapocho de Santiago)
# data rows; columns Fecha(D), Calidad(L), Volumen(M), Precio minimo(N),
# Precio maximo(O), Precio promedio ponderado(P), Origen(R) and Precio $/Kg(S)
# are reassigned row-by-row per the new ordering, while the descriptive
# columns (Mercado, Region, Codreg, Tipo, Producto, Categoria, Variedad,
# Unidad de comercializacion, Kg/unidad) stay put.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (was sourced from former row 34)
$ws.Range("D2").Value = 44294
$ws.Range("M2").Value = 480
$ws.Range("N2").Value = 7500
$ws.Range("O2").Value = 8000
$ws.Range("P2").Value = 7792
$ws.Range("R2").Value = 'Provincia de Linares'
$ws.Range("S2").Value = 3896

# Row 3  (was sourced from former row 20)
$ws.Range("D3").Value = 44165
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 25
$ws.Range("N3").Value = 10000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 10000
$ws.Range("S3").Value = 5000

# Row 4  (was sourced from former row 35)
$ws.Range("D4").Value = 44196
$ws.Range("M4").Value = 550
$ws.Range("N4").Value = 6500
$ws.Range("O4").Value = 7000
$ws.Range("P4").Value = 6818
$ws.Range("R4").Value = 'Provincia de Linares'
$ws.Range("S4").Value = 3409

# Row 5  (was sourced from former row 7)
$ws.Range("D5").Value = 44187
$ws.Range("M5").Value = 220
$ws.Range("N5").Value = 7000
$ws.Range("O5").Value = 7000
$ws.Range("P5").Value = 7000
$ws.Range("R5").Value = 'Provincia de Linares'
$ws.Range("S5").Value = 3500

# Row 6  (was sourced from former row 8)
$ws.Range("D6").Value = 44187
$ws.Range("L6").Value = 'Segunda'
$ws.Range("M6").Value = 260
$ws.Range("N6").Value = 5000
$ws.Range("O6").Value = 5000
$ws.Range("P6").Value = 5000
$ws.Range("R6").Value = 'Provincia de Linares'
$ws.Range("S6").Value = 2500

# Row 7  (was sourced from former row 16)
$ws.Range("D7").Value = 44365
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 10000
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 5000

# Row 8  (was sourced from former row 23)
$ws.Range("D8").Value = 44281
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 220
$ws.Range("N8").Value = 7500
$ws.Range("O8").Value = 7500
$ws.Range("P8").Value = 7500
$ws.Range("S8").Value = 3750

# Row 9  (was sourced from former row 15)
$ws.Range("D9").Value = 44364
$ws.Range("M9").Value = 75
$ws.Range("N9").Value = 10000
$ws.Range("O9").Value = 10000
$ws.Range("P9").Value = 10000
$ws.Range("R9").Value = 'Provincia de Curicó'
$ws.Range("S9").Value = 5000

# Row 10  (was sourced from former row 40)
$ws.Range("D10").Value = 44280
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 260
$ws.Range("N10").Value = 8000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 8000
$ws.Range("R10").Value = 'Provincia de Linares'
$ws.Range("S10").Value = 4000

# Row 11  (was sourced from former row 50)
$ws.Range("D11").Value = 44215
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 280
$ws.Range("N11").Value = 6600
$ws.Range("P11").Value = 6829
$ws.Range("S11").Value = 3414

# Row 12  (was sourced from former row 51)
$ws.Range("D12").Value = 44215
$ws.Range("L12").Value = 'Segunda'
$ws.Range("M12").Value = 560
$ws.Range("N12").Value = 5600
$ws.Range("O12").Value = 6000
$ws.Range("P12").Value = 5814
$ws.Range("R12").Value = 'Provincia de Linares'
$ws.Range("S12").Value = 2907

# Row 13  (was sourced from former row 30)
$ws.Range("D13").Value = 44188
$ws.Range("N13").Value = 6500
$ws.Range("O13").Value = 7000
$ws.Range("P13").Value = 6692
$ws.Range("S13").Value = 3346

# Row 14  (was sourced from former row 31)
$ws.Range("D14").Value = 44188
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 340
$ws.Range("N14").Value = 5000
$ws.Range("O14").Value = 5500
$ws.Range("P14").Value = 5206
$ws.Range("R14").Value = 'Provincia de Linares'
$ws.Range("S14").Value = 2603

# Row 15  (was sourced from former row 6)
$ws.Range("D15").Value = 44230
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 6000
$ws.Range("O15").Value = 6000
$ws.Range("P15").Value = 6000
$ws.Range("S15").Value = 3000

# Row 16  (was sourced from former row 22)
$ws.Range("D16").Value = 44224
$ws.Range("M16").Value = 420
$ws.Range("N16").Value = 6500
$ws.Range("O16").Value = 7000
$ws.Range("P16").Value = 6786
$ws.Range("R16").Value = 'Provincia de Linares'
$ws.Range("S16").Value = 3393

# Row 17  (was sourced from former row 52)
$ws.Range("D17").Value = 44286
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("R17").Value = 'Provincia de Linares'
$ws.Range("S17").Value = 4000

# Row 18  (was sourced from former row 24)
$ws.Range("D18").Value = 44358
$ws.Range("M18").Value = 50
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("R18").Value = 'Provincia de Curicó'
$ws.Range("S18").Value = 5000

# Row 19  (was sourced from former row 55)
$ws.Range("D19").Value = 44273
$ws.Range("M19").Value = 210
$ws.Range("N19").Value = 6000
$ws.Range("P19").Value = 6000
$ws.Range("S19").Value = 3000

# Row 20  (was sourced from former row 54)
$ws.Range("D20").Value = 44217
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 250
$ws.Range("N20").Value = 6500
$ws.Range("O20").Value = 6600
$ws.Range("P20").Value = 6560
$ws.Range("S20").Value = 3280

# Row 21  (was sourced from former row 37)
$ws.Range("D21").Value = 44357
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 10000
$ws.Range("P21").Value = 10000
$ws.Range("S21").Value = 5000

# Row 22  (was sourced from former row 29)
$ws.Range("D22").Value = 44202
$ws.Range("M22").Value = 310
$ws.Range("P22").Value = 6677
$ws.Range("S22").Value = 3338

# Row 23  (was sourced from former row 9)
$ws.Range("D23").Value = 44264
$ws.Range("N23").Value = 6000
$ws.Range("O23").Value = 6000
$ws.Range("P23").Value = 6000
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 3000

# Row 24  (was sourced from former row 10)
$ws.Range("D24").Value = 44174
$ws.Range("L24").Value = 'Especial'
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 7000
$ws.Range("O24").Value = 7000
$ws.Range("P24").Value = 7000
$ws.Range("S24").Value = 3500

# Row 25  (was sourced from former row 11)
$ws.Range("D25").Value = 44174
$ws.Range("L25").Value = 'Especial'
$ws.Range("M25").Value = 80
$ws.Range("N25").Value = 7000
$ws.Range("O25").Value = 7000
$ws.Range("P25").Value = 7000
$ws.Range("S25").Value = 3500

# Row 26  (was sourced from former row 41)
$ws.Range("D26").Value = 44293
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 8000
$ws.Range("O26").Value = 8000
$ws.Range("P26").Value = 8000
$ws.Range("S26").Value = 4000

# Row 27  (was sourced from former row 56)
$ws.Range("D27").Value = 44302
$ws.Range("M27").Value = 150
$ws.Range("N27").Value = 7000
$ws.Range("O27").Value = 7000
$ws.Range("P27").Value = 7000
$ws.Range("S27").Value = 3500

# Row 28  (was sourced from former row 47)
$ws.Range("D28").Value = 44292
$ws.Range("M28").Value = 120
$ws.Range("N28").Value = 8000
$ws.Range("O28").Value = 8000
$ws.Range("P28").Value = 8000
$ws.Range("S28").Value = 4000

# Row 29  (was sourced from former row 36)
$ws.Range("D29").Value = 44208
$ws.Range("M29").Value = 300
$ws.Range("N29").Value = 6000
$ws.Range("O29").Value = 6000
$ws.Range("P29").Value = 6000
$ws.Range("S29").Value = 3000

# Row 30  (was sourced from former row 42)
$ws.Range("D30").Value = 44204
$ws.Range("M30").Value = 120
$ws.Range("O30").Value = 6500
$ws.Range("P30").Value = 6500
$ws.Range("S30").Value = 3250

# Row 31  (was sourced from former row 43)
$ws.Range("D31").Value = 44204
$ws.Range("M31").Value = 150
$ws.Range("N31").Value = 7000
$ws.Range("O31").Value = 7000
$ws.Range("P31").Value = 7000
$ws.Range("S31").Value = 3500

# Row 32  (was sourced from former row 49)
$ws.Range("D32").Value = 44306
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 7000
$ws.Range("P32").Value = 7000
$ws.Range("R32").Value = 'Provincia de Curicó'
$ws.Range("S32").Value = 3500

# Row 33  (was sourced from former row 21)
$ws.Range("D33").Value = 44300
$ws.Range("M33").Value = 250
$ws.Range("N33").Value = 7000
$ws.Range("P33").Value = 7000
$ws.Range("R33").Value = 'Provincia de Curicó'
$ws.Range("S33").Value = 3500

# Row 34  (was sourced from former row 19)
$ws.Range("D34").Value = 44209
$ws.Range("M34").Value = 370
$ws.Range("N34").Value = 5800
$ws.Range("O34").Value = 6000
$ws.Range("P34").Value = 5935
$ws.Range("S34").Value = 2968

# Row 35  (was sourced from former row 17)
$ws.Range("D35").Value = 44237
$ws.Range("M35").Value = 150
$ws.Range("N35").Value = 6000
$ws.Range("O35").Value = 6000
$ws.Range("P35").Value = 6000
$ws.Range("R35").Value = 'Provincia de Colchagua'
$ws.Range("S35").Value = 3000

# Row 36  (was sourced from former row 18)
$ws.Range("D36").Value = 44237
$ws.Range("M36").Value = 200
$ws.Range("R36").Value = 'Región de O''Higgins'

# Row 37  (was sourced from former row 48)
$ws.Range("D37").Value = 44356
$ws.Range("M37").Value = 60

# Row 39  (was sourced from former row 32)
$ws.Range("D39").Value = 44216
$ws.Range("M39").Value = 340
$ws.Range("N39").Value = 6600
$ws.Range("O39").Value = 7000
$ws.Range("P39").Value = 6812
$ws.Range("S39").Value = 3406

# Row 40  (was sourced from former row 12)
$ws.Range("D40").Value = 44363
$ws.Range("M40").Value = 50
$ws.Range("N40").Value = 10000
$ws.Range("O40").Value = 10000
$ws.Range("P40").Value = 10000
$ws.Range("R40").Value = 'Provincia de Curicó'
$ws.Range("S40").Value = 5000

# Row 41  (was sourced from former row 14)
$ws.Range("D41").Value = 44299
$ws.Range("M41").Value = 250
$ws.Range("N41").Value = 7000
$ws.Range("O41").Value = 7000
$ws.Range("P41").Value = 7000
$ws.Range("R41").Value = 'Provincia de Curicó'
$ws.Range("S41").Value = 3500

# Row 42  (was sourced from former row 3)
$ws.Range("D42").Value = 44195
$ws.Range("M42").Value = 408
$ws.Range("N42").Value = 6509
$ws.Range("O42").Value = 7000
$ws.Range("P42").Value = 6774
$ws.Range("S42").Value = 3387

# Row 43  (was sourced from former row 5)
$ws.Range("D43").Value = 44239
$ws.Range("L43").Value = 'Primera'
$ws.Range("N43").Value = 6000
$ws.Range("O43").Value = 6000
$ws.Range("P43").Value = 6000
$ws.Range("R43").Value = 'Provincia de Curicó'
$ws.Range("S43").Value = 3000

# Row 44  (was sourced from former row 33)
$ws.Range("D44").Value = 44222
$ws.Range("L44").Value = 'Primera'
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 6500
$ws.Range("O44").Value = 7000
$ws.Range("P44").Value = 6800
$ws.Range("S44").Value = 3400

# Row 45  (was sourced from former row 2)
$ws.Range("D45").Value = 44362
$ws.Range("L45").Value = 'Primera'
$ws.Range("M45").Value = 75
$ws.Range("N45").Value = 10000
$ws.Range("O45").Value = 10000
$ws.Range("P45").Value = 10000
$ws.Range("R45").Value = 'Provincia de Curicó'
$ws.Range("S45").Value = 5000

# Row 46  (was sourced from former row 25)
$ws.Range("D46").Value = 44210
$ws.Range("M46").Value = 400
$ws.Range("N46").Value = 5800
$ws.Range("O46").Value = 6000
$ws.Range("P46").Value = 5910
$ws.Range("S46").Value = 2955

# Row 47  (was sourced from former row 53)
$ws.Range("D47").Value = 44301
$ws.Range("M47").Value = 200
$ws.Range("N47").Value = 7000
$ws.Range("O47").Value = 7000
$ws.Range("P47").Value = 7000
$ws.Range("R47").Value = 'Provincia de Curicó'
$ws.Range("S47").Value = 3500

# Row 48  (was sourced from former row 39)
$ws.Range("D48").Value = 44279
$ws.Range("M48").Value = 150
$ws.Range("N48").Value = 8000
$ws.Range("O48").Value = 8000
$ws.Range("P48").Value = 8000
$ws.Range("R48").Value = 'Provincia de Linares'
$ws.Range("S48").Value = 4000

# Row 49  (was sourced from former row 44)
$ws.Range("D49").Value = 44166
$ws.Range("L49").Value = 'Especial'
$ws.Range("M49").Value = 50
$ws.Range("N49").Value = 8000
$ws.Range("O49").Value = 8000
$ws.Range("P49").Value = 8000
$ws.Range("R49").Value = 'Provincia de Linares'
$ws.Range("S49").Value = 4000

# Row 50  (was sourced from former row 45)
$ws.Range("D50").Value = 44166
$ws.Range("L50").Value = 'Especial'
$ws.Range("M50").Value = 150
$ws.Range("N50").Value = 7200
$ws.Range("O50").Value = 7200
$ws.Range("P50").Value = 7200
$ws.Range("R50").Value = 'Región de O''Higgins'
$ws.Range("S50").Value = 3600

# Row 51  (was sourced from former row 46)
$ws.Range("D51").Value = 44166
$ws.Range("L51").Value = 'Primera'
$ws.Range("M51").Value = 80
$ws.Range("N51").Value = 7000
$ws.Range("O51").Value = 7000
$ws.Range("P51").Value = 7000
$ws.Range("S51").Value = 3500

# Row 52  (was sourced from former row 13)
$ws.Range("D52").Value = 44225
$ws.Range("M52").Value = 260
$ws.Range("N52").Value = 6000
$ws.Range("O52").Value = 6000
$ws.Range("P52").Value = 6000
$ws.Range("S52").Value = 3000

# Row 53  (was sourced from former row 27)
$ws.Range("D53").Value = 44238
$ws.Range("N53").Value = 6000
$ws.Range("O53").Value = 6000
$ws.Range("P53").Value = 6000
$ws.Range("S53").Value = 3000

# Row 54  (was sourced from former row 28)
$ws.Range("D54").Value = 44238
$ws.Range("M54").Value = 150
$ws.Range("N54").Value = 6000
$ws.Range("O54").Value = 6000
$ws.Range("P54").Value = 6000
$ws.Range("S54").Value = 3000

# Row 55  (was sourced from former row 4)
$ws.Range("D55").Value = 44236
$ws.Range("M55").Value = 450
$ws.Range("R55").Value = 'Provincia de Curicó'

# Row 56  (was sourced from former row 26)
$ws.Range("D56").Value = 44194
$ws.Range("M56").Value = 190
$ws.Range("N56").Value = 5800
$ws.Range("O56").Value = 6000
$ws.Range("P56").Value = 5916
$ws.Range("R56").Value = 'Provincia de Linares'
$ws.Range("S56").Value = 2958

Write-Host "Applied Frambuesa row permutation edit: 344 cells updated."
